# Applies the "Updated all sumsum.txt and excel files" edit:
#  - bump the "UPDATED" footer date string
#  - change the highlight fill color used on the footer row
#  - update the raw benchmark numbers on rows 16-24 (columns E/H/J/L);
#    the dependent formulas (F16:F24 and the summary table in rows 32-43)
#    recalculate automatically
#  - add an (empty, but styled) helper column M across rows 15-25
#  - move the active selection to M15:M25

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- footer date string -------------------------------------------------
$ws.Cells.Item(44, 1).Value = "UPDATED 20180326"

# --- footer highlight color (was FF6600, now FF0000) --------------------
# OLE colors are 0x00BBGGRR, so pure red (FF0000) is 255.
$ws.Range("A44:B44").Interior.Color = 255

# --- updated raw data, rows 16-24 (E, H, J, L columns) -------------------
$ws.Range("E16").Value = 773.8
$ws.Range("H16").Value = 360.5
$ws.Range("J16").Value = 289
$ws.Range("L16").Value = 774.4

$ws.Range("E17").Value = 363.2
$ws.Range("H17").Value = 231.4
$ws.Range("J17").Value = 75.9
$ws.Range("L17").Value = 363.4

$ws.Range("E18").Value = 307.7
$ws.Range("H18").Value = 194.1
$ws.Range("J18").Value = 58.5
$ws.Range("L18").Value = 307.9

$ws.Range("E19").Value = 203
$ws.Range("H19").Value = 89.5
$ws.Range("J19").Value = 58.5
$ws.Range("L19").Value = 203.2

$ws.Range("E20").Value = 195.5
$ws.Range("H20").Value = 81.9
$ws.Range("J20").Value = 58.5
$ws.Range("L20").Value = 195.7

$ws.Range("E21").Value = 175
$ws.Range("H21").Value = 63.9
$ws.Range("J21").Value = 56.1
$ws.Range("L21").Value = 175.2

$ws.Range("E22").Value = 167.7
$ws.Range("H22").Value = 63.4
$ws.Range("J22").Value = 52.9
$ws.Range("L22").Value = 167.9

$ws.Range("E23").Value = 157.9
$ws.Range("H23").Value = 63.4
$ws.Range("J23").Value = 42.3
$ws.Range("L23").Value = 158

$ws.Range("E24").Value = 135.1
$ws.Range("H24").Value = 63.4
$ws.Range("J24").Value = 42.3
$ws.Range("L24").Value = 135.3

# --- new helper column M, rows 15-25 (empty cells, "no fill" style) ------
$ws.Range("M15:M25").Interior.ColorIndex = -4142

# --- move the selection to match the author's final cursor position ------
$ws.Range("M15:M25").Select()
